$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Emmett"
$ws.Range("B2").Value = "Pagac"
$ws.Range("C2").Value = "ryan.lang@gmail.com"
$ws.Range("D2").Value = "f58fwlrf5c"
